$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
Write-Host "ScrollRow=" $win.ScrollRow
Write-Host "ScrollColumn=" $win.ScrollColumn
$win.ScrollRow = 99
$win.ScrollColumn = 9
Write-Host "ScrollRow2=" $win.ScrollRow
Write-Host "ScrollColumn2=" $win.ScrollColumn
